# Regenerated save_data: column G ("K") is recomputed (now derived from
# actual strikeouts "K" instead of the old "Strike#" count) for every
# data row (rows 2-59) on the active sheet. Write the recalculated
# per-row K values in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 2
    3 = 1
    4 = 4
    5 = 2
    6 = 0
    7 = 3
    8 = 1
    9 = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 1
    17 = 5
    18 = 0
    19 = 0
    20 = 3
    21 = 2
    22 = 2
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 3
    40 = 0
    41 = 2
    42 = 0
    43 = 3
    44 = 1
    45 = 2
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 2
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 1
    57 = 1
    58 = 3
    59 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
